# Append new daily snapshot rows 58-63 to the Positions sheet (sheet1),
# mirroring the existing per-column layout. Extends the used range from
# A1:Z57 to A1:Z63. Blank/omitted columns in a given day are left empty,
# matching how earlier rows (e.g. row 57) already skip zero-valued symbols.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A on existing rows uses a date-formatted style (s="2"). Copy that
# style from the last existing row (A57) onto the new A-column cells before
# writing their values, so formatting matches exactly.
$ws.Range("A57").Copy()
$ws.Range("A58:A63").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
    @{ Row = 58; Cells = @{ "A" = 45488; "B" = 690.8162800136; "C" = 241.5255482265; "D" = 0; "E" = 0; "F" = 0; "G" = 131.5124349; "H" = 0; "I" = 233.7179893178; "J" = 0; "K" = 21.457450160069; "L" = 0; "M" = 0; "N" = 170.1861951616; "O" = 59.45206406099999; "P" = 0; "Q" = 0.0000027168; "R" = 0; "S" = 0; "T" = 0; "U" = 402.5926200617187; "V" = 0; "W" = 0; "X" = 0; "Y" = 0; "Z" = 257.207866961948 } },
    @{ Row = 59; Cells = @{ "A" = 45489; "C" = 238.8034031255; "D" = 0; "E" = 0; "G" = 136.59134985; "H" = 0; "J" = 0; "K" = 20.945948703264; "L" = 0; "M" = 0; "N" = 158.39376589056; "P" = 0; "Q" = 0.0000029688; "S" = 0; "T" = 0; "X" = 0; "Z" = 270.503112633254 } },
    @{ Row = 60; Cells = @{ "A" = 45490; "B" = 684.0264983876; "C" = 234.8456842675; "D" = 0; "E" = 0; "F" = 0; "G" = 128.41115055; "H" = 0; "I" = 228.040664475; "J" = 0; "K" = 20.8782794757346; "L" = 0; "M" = 0; "N" = 154.80913540192; "O" = 57.786227993; "P" = 0; "Q" = 0.0000028056; "R" = 0; "S" = 0; "T" = 0; "U" = 370.4824365105616; "V" = 0; "W" = 0; "X" = 0; "Y" = 0; "Z" = 269.275539089268 } },
    @{ Row = 61; Cells = @{ "A" = 45491; "B" = 682.9584272608; "C" = 237.581003275; "D" = 0; "E" = 0; "G" = 126.5683584; "H" = 0; "I" = 233.5023825694; "J" = 0; "K" = 21.0216163379646; "L" = 0; "M" = 0; "N" = 162.38041101344; "O" = 58.09095410299999; "P" = 0; "Q" = 0.0000027936; "R" = 0; "S" = 0; "T" = 0; "W" = 0; "X" = 0; "Y" = 0; "Z" = 269.75408470811 } },
    @{ Row = 62; Cells = @{ "A" = 45492; "B" = 711.4781784; "C" = 242.9219823155; "D" = 0; "E" = 0; "F" = 0; "G" = 127.69201215; "H" = 0; "I" = 248.0425214917; "J" = 0; "K" = 21.0492955163586; "L" = 0; "M" = 0; "N" = 161.44237686688; "O" = 60.274824558; "P" = 0; "Q" = 0.0000030456; "R" = 0; "S" = 0; "T" = 0; "U" = 377.9023195622234; "V" = 0; "W" = 0; "X" = 0; "Y" = 0; "Z" = 266.799585670042 } },
    @{ Row = 63; Cells = @{ "A" = 45493; "B" = 716.6009066704001; "C" = 243.890611125; "D" = 0; "E" = 0; "G" = 126.52341225; "H" = 0; "I" = 254.605042776; "J" = 0; "K" = 21.0535290986384; "L" = 0; "M" = 0; "N" = 160.03532564704; "O" = 60.244351947; "P" = 0; "Q" = 0.0000029328; "R" = 0; "S" = 0; "T" = 0; "W" = 0; "X" = 0; "Y" = 0; "Z" = 277.639684253376 } }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    foreach ($col in $r.Cells.Keys) {
        $ws.Range("$col$rowNum").Value = $r.Cells[$col]
    }
}
